$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("continent")

$ws.Range("A8").Value = "EU_OECD"
$ws.Range("B8").Value = 516
$ws.Range("C8").Value = 388
$ws.Range("D8").Value = 516
$ws.Range("E8").Value = 0.751937984496124
$ws.Range("F8").Value = 1
